$wb = $excel.ActiveWorkbook

# School sheet: update B1 value from 200 to 250
$school = $wb.Worksheets.Item("School")
$school.Range("B1").Value = 250

# Purchase sheet: update C2 and C4 values, and move selection to C4
$purchase = $wb.Worksheets.Item("Purchase")
$purchase.Range("C2").Value = 16
$purchase.Range("C4").Value = 10
$purchase.Range("C4").Select()

# Tool sheet: update B1 and B5 values, and move selection to B1
$tool = $wb.Worksheets.Item("Tool")
$tool.Range("B1").Value = 0
$tool.Range("B5").Value = 62

# Make Tool the active sheet (this also updates tabSelected/activeTab and selection)
$tool.Activate()
$tool.Range("B1").Select()
